$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.595.95"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "3.815.39"
$ws.Range("E3").Value = "  +0.35%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "706.50"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.36%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "174.95"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("D7").Value = "3.814.86"
$ws.Range("E7").Value = "  +0.37%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +2.46%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.31"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.56%  "
$ws.Range("E12").Value = "  +0.48%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000261"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +7.66%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "36.47"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "4.458.56"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "3.836.61"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "71.552.32"
$ws.Range("E17").Value = "  +2.13%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "17.77"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.22"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("E20").Value = "  +0.17%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.93"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.86%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "484.13"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "84.62"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("E25").Value = "  -0.93%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.36"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "3.966.44"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  +12.57%  "
$ws.Range("E31").Value = "  +0.04%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.64"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.85%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.31"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  +7.01%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "29.67"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +1.70%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.47"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.47%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.04"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.01%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +11.13%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.989"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.31%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.000315"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +16.03%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "164.82"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.68%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "45.02"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "48.73"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.68%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "419.91"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +7.77%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  -2.11%  "
